# Product Backlog update for this sprint:
# refresh the "Priority" column (A) for a handful of backlog items and
# leave the cursor where the team left off reviewing the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A7").Value = 0
$ws.Range("A8").Value = 0
$ws.Range("A12").Value = 1
$ws.Range("A15").Value = 1

$ws.Range("A11").Select()
